
$d = $word.ActiveDocument


$r = $d.Content
$r.Find.Execute("Regras_de_Negócio", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $r.Paragraphs(1)
$para.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00842E3F" w:rsidRPr="00842E3F" w:rsidRDefault="00920E7B" w:rsidP="00842E3F"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Regras_de_Negócio</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')


$r = $d.Content
$r.Find.Execute("Consultar Livros, Solicitar Empréstimo, Realizar Devolução", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $r.Paragraphs(1)
$para.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00670C58" w:rsidRDefault="00670C58" w:rsidP="00A334D8"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="12"/></w:numPr><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Aluno:</w:t></w:r><w:r><w:t xml:space="preserve"> Consultar Livros, Solicitar Empréstimo, Realizar Devolução</w:t></w:r><w:r><w:t>, Informar Pagamento</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>')


$r = $d.Content
$r.Find.Execute("Um aluno, cuja situação da matrícula for diferente de ativa", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $r.Paragraphs(1)
$para.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00920E7B" w:rsidRDefault="00920E7B" w:rsidP="00920E7B"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:tab/><w:t>Um aluno, cuja situação da matrícula for diferente de ativa, não poderá efetuar empréstimos.</w:t></w:r></w:p>')


$r = $d.Content
$r.Find.Execute("endereço deve ser persistido no sistema", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $r.Paragraphs(1)
$para.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00240F6C" w:rsidRDefault="00240F6C" w:rsidP="00240F6C"><w:pPr><w:ind w:firstLine="708"/><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">O </w:t></w:r><w:r w:rsidRPr="00240F6C"><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">endereço deve ser persistido no sistema, porém os dados devem ser obtidos através de consulta ao </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00240F6C"><w:rPr><w:i/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>WebServices</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00240F6C"><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> dos Correios. Caso o endereço não exista nos correios, o sistema deve permitir o cadastramento do mesmo.</w:t></w:r></w:p>')


$r = $d.Content
$r.Find.Execute("O sistema deve verificar se não existe", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para1 = $r.Paragraphs(1)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$wholeRange = $d.Range($para1.Range.Start, $lastPara.Range.End)
$wholeRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00332A64" w:rsidRPr="00332A64" w:rsidRDefault="00332A64" w:rsidP="00332A64"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">O sistema deve verificar se não existe nenhum vínculo de empréstimo com o livro no momento de exclusão, caso exista o sistema não pode permitir a exclusão, exibindo a mensagem </w:t></w:r><w:r w:rsidR="00DB1989"><w:rPr><w:b/></w:rPr><w:t>MSG10</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>RN11 – REALIZAÇÃO DE EMPRÉSTIMO</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="708"/><w:jc w:val="both"/><w:rPr><w:rFonts w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="24292E"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Para cada empréstimo devem ser registradas as seguintes informações: data do </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="24292E"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>empréstimo</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="24292E"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">, data prevista para devolução (30 dias após a data do </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="24292E"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">empréstimo), data da devolução efetiva, para qual Aluno e qual </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="24292E"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>funcionário</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:color w:val="24292E"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> registrou este empréstimo).</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00332A64" w:rsidRPr="006D7815" w:rsidRDefault="00332A64" w:rsidP="006D7815"><w:bookmarkStart w:id="11" w:name="_GoBack"/><w:bookmarkEnd w:id="11"/></w:p>')


$h = $d.Sections(1).Headers(1)
$r = $h.Range
$r.Find.Execute("Regras_de_Negócio", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $r.Paragraphs(1)
$para.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00842E3F" w:rsidRPr="00E238ED" w:rsidRDefault="00842E3F" w:rsidP="00920E7B"><w:pPr><w:tabs><w:tab w:val="left" w:pos="1135"/></w:tabs><w:spacing w:before="60" w:after="60"/><w:ind w:right="68"/><w:rPr><w:rFonts w:cs="Arial"/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="00E238ED"><w:rPr><w:rFonts w:cs="Arial"/><w:bCs/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidRPr="00E238ED"><w:rPr><w:rFonts w:cs="Arial"/><w:bCs/></w:rPr><w:instrText xml:space="preserve"> TITLE   \* MERGEFORMAT </w:instrText></w:r><w:r w:rsidRPr="00E238ED"><w:rPr><w:rFonts w:cs="Arial"/><w:bCs/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:bCs/></w:rPr><w:t>SGB_</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00920E7B"><w:rPr><w:rFonts w:cs="Arial"/><w:bCs/></w:rPr><w:t>Regras_de_Negócio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00E238ED"><w:rPr><w:rFonts w:cs="Arial"/><w:bCs/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>')
